$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (avoid Excel auto-converting numeric-looking strings
# like "1.002" or "0.000009932" into real numbers, which would reformat them,
# e.g. into scientific notation) while writing the new values, then restore
# the default "Normal" style so the cells end up unstyled, as in the source file.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.866.71"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "1.879.08"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "324.44"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "0.4617"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").Value = "0.3879"
$ws.Range("E8").Value = "  -1.93%  "
$ws.Range("D9").Value = "0.07852"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "0.9833"
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").Value = "21.76"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "1.866.61"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("D13").Value = "7.003"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "5.672"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "0.06976"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "88.61"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").Value = "0.000009932"
$ws.Range("E18").Value = "  -2.62%  "
$ws.Range("D19").Value = "16.95"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "28.882.68"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "5.272"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").Value = "10.99"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("D24").Value = "2.103"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("D25").Value = "155.82"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "19.38"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").Value = "5.900"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "117.83"
$ws.Range("E28").Value = "  -2.48%  "
$ws.Range("D29").Value = "1.900"
$ws.Range("E29").Value = "  -6.51%  "
$ws.Range("D30").Value = "0.09366"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("D32").Value = "5.269"
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("D33").Value = "1.319"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").Value = "3.249"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "0.05752"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").Value = "1.170"
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("D37").Value = "0.02077"
$ws.Range("E37").Value = "  -1.30%  "
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("D39").Value = "7.633"
$ws.Range("E39").Value = "  -6.45%  "
$ws.Range("D40").Value = "0.5664"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").Value = "0.1775"
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("D42").Value = "9.684"
$ws.Range("E42").Value = "  -3.79%  "
$ws.Range("D43").Value = "11.92"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").Value = "2.227"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "0.5333"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").Value = "0.07040"
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "1.845"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "2.547"
$ws.Range("E48").Value = "  +1.76%  "
$ws.Range("D49").Value = "112.38"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "1.067"
$ws.Range("E50").Value = "  -5.75%  "
$ws.Range("D51").Value = "70.94"
$ws.Range("E51").Value = "  -1.56%  "

$ws.Range("D2:E51").Style = "Normal"
